# ProvarCache data update
# Updates the "RMA Details Maintenance Grid" sheet so the cached RMA
# receipt-line values point at the newly created RMA "RMA-D2PI" (replacing
# the previous run's "RMA-52JG" values) for each of the 3 maintenance rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RMA Details Maintenance Grid")

# Row 2 (RMA line 1)
$ws.Range("E2").Value = "RMA-D2PI-001"
$ws.Range("F2").Value = "RMA-D2PI-1-1"
$ws.Range("J2").Value = "a7s5f000000xKLvAAM"

# Row 3 (RMA line 2)
$ws.Range("E3").Value = "RMA-D2PI-002"
$ws.Range("F3").Value = "RMA-D2PI-1-2"
$ws.Range("J3").Value = "a7s5f000000xKLwAAM"

# Row 4 (RMA line 3)
$ws.Range("E4").Value = "RMA-D2PI-003"
$ws.Range("F4").Value = "RMA-D2PI-1-3"
$ws.Range("J4").Value = "a7s5f000000xKLxAAM"
